$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Cells changing between numeric and text representation ---
# (copy number-format from a stable same-style cell, then set the new value)
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = "0"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "0"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "***.*"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1

$ws.Range("D15").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = "0"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value = "0"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Value = "***.*"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1

$ws.Range("K14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = 0

$ws.Range("I14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Value = 1

$ws.Range("I14").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = 1

$ws.Range("I14").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Value = 1

$ws.Range("I14").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("F30").Value = 1

$ws.Range("D15").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").Value = "0"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "***.*"

# --- Simple value-only changes (style/type unchanged) ---
# Row 15
$ws.Range("M15").Value = 112.5

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 125
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 27.777777777777
$ws.Range("I16").Value = 172
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = 3.614457831325
$ws.Range("L16").Value = 11.688311688311
$ws.Range("M16").Value = -48.036253776435
$ws.Range("N16").Value = -84.587813620071

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 270
$ws.Range("J17").Value = 265
$ws.Range("K17").Value = 1.88679245283
$ws.Range("L17").Value = 8.433734939759
$ws.Range("M17").Value = 52.542372881355
$ws.Range("N17").Value = -47.674418604651

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = 262
$ws.Range("J18").Value = 243
$ws.Range("K18").Value = 7.818930041152
$ws.Range("L18").Value = -1.872659176029
$ws.Range("M18").Value = -35.941320293398
$ws.Range("N18").Value = -77.721088435374

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = -44.285714285714
$ws.Range("I19").Value = 616
$ws.Range("J19").Value = 672
$ws.Range("K19").Value = -8.333333333333
$ws.Range("L19").Value = 4.054054054054
$ws.Range("M19").Value = 43.255813953488
$ws.Range("N19").Value = 31.343283582089

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 22.222222222222
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 154
$ws.Range("K20").Value = -31.168831168831
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -23.741007194244
$ws.Range("N20").Value = -86.733416770963

# Row 21
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 19.230769230769
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -22.463768115942
$ws.Range("I21").Value = 1446
$ws.Range("J21").Value = 1516
$ws.Range("K21").Value = -4.617414248021
$ws.Range("L21").Value = 0.696378830083
$ws.Range("M21").Value = -3.342245989304
$ws.Range("N21").Value = -64.987893462469

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = -9.090909090909
$ws.Range("L22").Value = -44.444444444444
$ws.Range("M22").Value = -48.717948717948

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 147
$ws.Range("J23").Value = 171
$ws.Range("K23").Value = -14.035087719298
$ws.Range("L23").Value = 3.521126760563
$ws.Range("M23").Value = 22.5

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 77.777777777777
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 36.842105263157
$ws.Range("I24").Value = 952
$ws.Range("J24").Value = 897
$ws.Range("K24").Value = 6.13154960981
$ws.Range("L24").Value = -9.333333333333
$ws.Range("M24").Value = -11.028037383177

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 500
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 76.923076923076
$ws.Range("I25").Value = 287
$ws.Range("J25").Value = 137
$ws.Range("K25").Value = 109.489051094891
$ws.Range("L25").Value = 21.097046413502

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -7.142857142857
$ws.Range("I26").Value = 451
$ws.Range("J26").Value = 430
$ws.Range("K26").Value = 4.883720930232
$ws.Range("L26").Value = 3.203661327231
$ws.Range("M26").Value = 7.637231503579

# Row 27
$ws.Range("F27").Value = 2

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 60
$ws.Range("I28").Value = 65
$ws.Range("J28").Value = 49
$ws.Range("K28").Value = 32.653061224489
$ws.Range("L28").Value = 47.727272727272

# Row 29
$ws.Range("I29").Value = 9
$ws.Range("K29").Value = -25
$ws.Range("L29").Value = -10
$ws.Range("M29").Value = -40
$ws.Range("N29").Value = -88.607594936708

# Row 30
$ws.Range("I30").Value = 8
$ws.Range("K30").Value = -11.111111111111
$ws.Range("L30").Value = -20
$ws.Range("M30").Value = -42.857142857142
$ws.Range("N30").Value = -89.041095890411

# Row 31
$ws.Range("L31").Value = -27.586206896551
